$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear existing contents for the data region that will be rewritten (rows 20-42)
$ws.Range("A20:G42").ClearContents()

# Rewrite rows 20-42 with the de-fogged, re-dated data
$ws.Range("A20").Value = 42552
$ws.Range("B20").Value = 629
$ws.Range("C20").Value = 27
$ws.Range("D20").Value = 115
$ws.Range("E20").Value = 12
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = "3 ATPU, 1 seal, 1 COLO"

$ws.Range("A21").Value = 42557
$ws.Range("B21").Value = 529
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 71
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 15
$ws.Range("G21").Value = "2 ATPU"

$ws.Range("A22").Value = 42558
$ws.Range("B22").Value = 379
$ws.Range("C22").Value = 60
$ws.Range("D22").Value = 82
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = "2 ATPU, 1 DCCO"

$ws.Range("A23").Value = 42561
$ws.Range("B23").Value = 506
$ws.Range("C23").Value = 39
$ws.Range("D23").Value = 141
$ws.Range("E23").Value = 35
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = "1 canada goose,3 DCCO, 1 gray seal, 1 COMU, 5 ATPU"

$ws.Range("A24").Value = 42562
$ws.Range("B24").Value = 552
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 75
$ws.Range("E24").Value = 53
$ws.Range("F24").Value = 19
$ws.Range("G24").Value = "fog"

$ws.Range("A25").Value = 42563
$ws.Range("B25").Value = 532
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 113
$ws.Range("E25").Value = 29
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = "1 COMU, 23 COTE"

$ws.Range("A26").Value = 42564
$ws.Range("B26").Value = 545
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 124
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = "1 DCC0"

$ws.Range("A27").Value = 42567
$ws.Range("B27").Value = 550
$ws.Range("C27").Value = 51
$ws.Range("D27").Value = 122
$ws.Range("E27").Value = 19
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = "1 ATPU, 1 seal, 1 DCCO, 1 wimbrel"

$ws.Range("A28").Value = 42568
$ws.Range("B28").Value = 589
$ws.Range("C28").Value = 43
$ws.Range("D28").Value = 112
$ws.Range("E28").Value = 43
$ws.Range("F28").Value = 8

$ws.Range("A29").Value = 42570
$ws.Range("B29").Value = 520
$ws.Range("C29").Value = 76
$ws.Range("D29").Value = 115
$ws.Range("E29").Value = 34
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = "3 DCCO,2  ATPU, 1 gray seal"

$ws.Range("A30").Value = 42571
$ws.Range("B30").Value = 352
$ws.Range("C30").Value = 27
$ws.Range("D30").Value = 70
$ws.Range("E30").Value = 72
$ws.Range("G30").Value = "stop distinguishing chicks; 3 red knots, 3 DCCO"

$ws.Range("A31").Value = 42572
$ws.Range("B31").Value = 482
$ws.Range("C31").Value = 34
$ws.Range("D31").Value = 86
$ws.Range("E31").Value = 30
$ws.Range("G31").Value = "6 DCCO, 1 ATPU"

$ws.Range("A32").Value = 42573
$ws.Range("B32").Value = 435
$ws.Range("C32").Value = 32
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = 32
$ws.Range("G32").Value = "5 DCCO, 1 peregrine"

$ws.Range("A33").Value = 42574
$ws.Range("B33").Value = 470
$ws.Range("C33").Value = 47
$ws.Range("D33").Value = 110
$ws.Range("E33").Value = 54
$ws.Range("G33").Value = "1 DCCO"

$ws.Range("A34").Value = 42575
$ws.Range("B34").Value = 476
$ws.Range("C34").Value = 34
$ws.Range("D34").Value = 77
$ws.Range("E34").Value = 40
$ws.Range("G34").Value = "1 DCCO, 1 seal, 1 ATPU"

$ws.Range("A35").Value = 42576
$ws.Range("B35").Value = 535
$ws.Range("C35").Value = 31
$ws.Range("D35").Value = 177
$ws.Range("E35").Value = 52
$ws.Range("G35").Value = "1 seal, 1 ATPU, 1 COLO"

$ws.Range("A36").Value = 42578
$ws.Range("B36").Value = 581
$ws.Range("C36").Value = 65
$ws.Range("D36").Value = 161
$ws.Range("E36").Value = 34
$ws.Range("G36").Value = "1 ATPU, 1 LAGU, 1 seal, 2 least sandpipers"

$ws.Range("A37").Value = 42579
$ws.Range("B37").Value = 633
$ws.Range("C37").Value = 58
$ws.Range("D37").Value = 166
$ws.Range("E37").Value = 51
$ws.Range("G37").Value = "1 ATPU, 3 BAEA, 1 DCCO"

$ws.Range("A38").Value = 42581
$ws.Range("B38").Value = 625
$ws.Range("C38").Value = 111
$ws.Range("D38").Value = 181
$ws.Range("E38").Value = 75
$ws.Range("G38").Value = "3 BAEA, 1 DCCO"

$ws.Range("A39").Value = 42582
$ws.Range("B39").Value = 625
$ws.Range("C39").Value = 39
$ws.Range("D39").Value = 241
$ws.Range("E39").Value = 68
$ws.Range("G39").Value = "1 DCCO, 1 ATPU, 1 seal"

$ws.Range("A40").Value = 42583
$ws.Range("B40").Value = 608
$ws.Range("C40").Value = 40
$ws.Range("D40").Value = 274
$ws.Range("E40").Value = 67
$ws.Range("G40").Value = "1 ATPU, 1 DCCO"

$ws.Range("A41").Value = 42584
$ws.Range("B41").Value = 594
$ws.Range("C41").Value = 39
$ws.Range("D41").Value = 159
$ws.Range("E41").Value = 49
$ws.Range("G41").Value = "2 DCCO, 1 ATPU, 1 BAEA"

$ws.Range("A42").Value = 42585
$ws.Range("B42").Value = 620
$ws.Range("C42").Value = 30
$ws.Range("D42").Value = 325
$ws.Range("E42").Value = 84
$ws.Range("G42").Value = "count is post-eagle; 1 ATPU, 1 DCCO, 2 BAEA, 4 COLO, 1 seal"

# Remove now-obsolete trailing rows (old rows 43-55)
$ws.Rows("43:55").Delete()

# Update the sheet view to match the post-edit scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A38:XFD38").Select()